# Remove the block-diagram image (inline picture) from the document.
# The picture was an inline <w:drawing> run inside a paragraph that also
# contains the "_GoBack" bookmark; the paragraph itself (and the bookmark)
# must remain intact, only the picture run is deleted.

$d = $word.ActiveDocument

if ($d.InlineShapes.Count -gt 0) {
    for ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {
        $d.InlineShapes.Item($i).Delete()
    }
}
